# Atualização automática de VIAMAO.xlsx
#
# - Remove the "Desarquivamentos Pendentes" sheet (no longer needed).
# - Rename "Paineis DARQ" -> "PAINEIS DARQ" (uppercase tab name).
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO".

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Drop the obsolete "Desarquivamentos Pendentes" worksheet entirely.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Rename the remaining sheets that changed their tab names.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
